$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Agus Priyanto"
$ws.Range("D4").Value = "Agus Priyanto"
$ws.Range("D5").Value = "Agus Priyanto"
$ws.Range("D6").Value = "Agus Priyanto"
$ws.Range("D7").Value = "Agus Priyanto"

$ws.Range("L16").Select()
